# ====================================================================
# Costa Rica Primera Division - league base update (07-03-2024 23:43)
# ====================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that had their home/away match data swapped (source rows reordered) ---
# Row 38 <-> Row 39
$ws.Range("B38").Value = 6782522
$ws.Range("F38").Value = "Municipal Perez Zeledon"
$ws.Range("G38").Value = "Sporting San Jose"
$ws.Range("I38").Value = 2
$ws.Range("J38").Value = "A"
$ws.Range("K38").Value = 2.5
$ws.Range("L38").Value = 3.5
$ws.Range("M38").Value = 2.5
$ws.Range("N38").Value = 2.2
$ws.Range("O38").Value = 3.5
$ws.Range("P38").Value = 2.9
$ws.Range("R38").Value = 1.9
$ws.Range("S38").Value = 1.9
$ws.Range("T38").Value = 2.5
$ws.Range("W38").Value = -1
$ws.Range("Y38").Value = 1.9
$ws.Range("Z38").Value = -1
$ws.Range("AA38").Value = 0.8999999999999999
$ws.Range("AB38").Value = 0.8999999999999999
$ws.Range("AC38").Value = -1
$ws.Range("B39").Value = 6781354
$ws.Range("F39").Value = "Puntarenas"
$ws.Range("G39").Value = "AD San Carlos"
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = "H"
$ws.Range("K39").Value = 2.4
$ws.Range("L39").Value = 3.2
$ws.Range("M39").Value = 2.8
$ws.Range("N39").Value = 2.3
$ws.Range("O39").Value = 3.2
$ws.Range("P39").Value = 3
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 1.8
$ws.Range("T39").Value = 2.25
$ws.Range("W39").Value = 1.3
$ws.Range("Y39").Value = -1
$ws.Range("Z39").Value = 1
$ws.Range("AA39").Value = -1
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 0.8999999999999999

# Row 91 <-> Row 92
$ws.Range("B91").Value = 6782566
$ws.Range("F91").Value = "Cartagines"
$ws.Range("G91").Value = "Deportivo Saprissa"
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 4
$ws.Range("J91").Value = "A"
$ws.Range("K91").Value = 3.2
$ws.Range("L91").Value = 3.4
$ws.Range("M91").Value = 2
$ws.Range("N91").Value = 2.9
$ws.Range("O91").Value = 3.5
$ws.Range("P91").Value = 2.15
$ws.Range("Q91").Value = 0.25
$ws.Range("R91").Value = 1.875
$ws.Range("S91").Value = 1.925
$ws.Range("T91").Value = 3
$ws.Range("U91").Value = 1.975
$ws.Range("V91").Value = 1.825
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = 1.15
$ws.Range("AA91").Value = 0.925
$ws.Range("AB91").Value = 0.9750000000000001
$ws.Range("AC91").Value = -1
$ws.Range("B92").Value = 6782568
$ws.Range("F92").Value = "Sporting San Jose"
$ws.Range("G92").Value = "AD Guanacasteca"
$ws.Range("H92").Value = 1
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = "D"
$ws.Range("K92").Value = 1.909
$ws.Range("L92").Value = 3.6
$ws.Range("M92").Value = 3.3
$ws.Range("N92").Value = 2
$ws.Range("O92").Value = 3.6
$ws.Range("P92").Value = 3.1
$ws.Range("Q92").Value = -0.5
$ws.Range("R92").Value = 2
$ws.Range("S92").Value = 1.8
$ws.Range("T92").Value = 2.5
$ws.Range("U92").Value = 1.825
$ws.Range("V92").Value = 1.975
$ws.Range("X92").Value = 2.6
$ws.Range("Y92").Value = -1
$ws.Range("AA92").Value = 0.8
$ws.Range("AB92").Value = -1
$ws.Range("AC92").Value = 0.9750000000000001

# Row 95 <-> Row 96
$ws.Range("B95").Value = 6782565
$ws.Range("F95").Value = "Santos de Gupiles"
$ws.Range("G95").Value = "Municipal Perez Zeledon"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = "H"
$ws.Range("K95").Value = 1.833
$ws.Range("L95").Value = 3.4
$ws.Range("M95").Value = 3.75
$ws.Range("N95").Value = 1.833
$ws.Range("P95").Value = 3.75
$ws.Range("Q95").Value = -0.5
$ws.Range("R95").Value = 1.875
$ws.Range("S95").Value = 1.925
$ws.Range("U95").Value = 2
$ws.Range("V95").Value = 1.8
$ws.Range("W95").Value = 0.833
$ws.Range("Y95").Value = -1
$ws.Range("Z95").Value = 0.875
$ws.Range("AA95").Value = -1
$ws.Range("AB95").Value = -1
$ws.Range("AC95").Value = 0.8
$ws.Range("B96").Value = 6782567
$ws.Range("F96").Value = "AD Grecia"
$ws.Range("G96").Value = "Municipal Liberia"
$ws.Range("I96").Value = 3
$ws.Range("J96").Value = "A"
$ws.Range("K96").Value = 2.875
$ws.Range("L96").Value = 3.5
$ws.Range("M96").Value = 2.15
$ws.Range("N96").Value = 2.3
$ws.Range("P96").Value = 2.6
$ws.Range("Q96").Value = 0
$ws.Range("R96").Value = 1.8
$ws.Range("S96").Value = 2
$ws.Range("U96").Value = 1.8
$ws.Range("V96").Value = 2
$ws.Range("W96").Value = -1
$ws.Range("Y96").Value = 1.6
$ws.Range("Z96").Value = -1
$ws.Range("AA96").Value = 1
$ws.Range("AB96").Value = 0.8
$ws.Range("AC96").Value = -1

# Row 130 <-> Row 131
$ws.Range("B130").Value = 6782595
$ws.Range("F130").Value = "Herediano"
$ws.Range("G130").Value = "Sporting San Jose"
$ws.Range("H130").Value = 3
$ws.Range("K130").Value = 1.4
$ws.Range("L130").Value = 4.75
$ws.Range("M130").Value = 7
$ws.Range("N130").Value = 1.363
$ws.Range("O130").Value = 4.75
$ws.Range("P130").Value = 8.5
$ws.Range("Q130").Value = -1.25
$ws.Range("T130").Value = 3
$ws.Range("U130").Value = 1.95
$ws.Range("V130").Value = 1.85
$ws.Range("W130").Value = 0.363
$ws.Range("AB130").Value = 0
$ws.Range("AC130").Value = -0
$ws.Range("B131").Value = 6782598
$ws.Range("F131").Value = "Municipal Perez Zeledon"
$ws.Range("G131").Value = "Cartagines"
$ws.Range("H131").Value = 1
$ws.Range("K131").Value = 4.5
$ws.Range("L131").Value = 3.75
$ws.Range("M131").Value = 1.615
$ws.Range("N131").Value = 3.4
$ws.Range("O131").Value = 3.4
$ws.Range("P131").Value = 1.85
$ws.Range("Q131").Value = 0.5
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 1.9
$ws.Range("V131").Value = 1.9
$ws.Range("W131").Value = 2.4
$ws.Range("AB131").Value = -1
$ws.Range("AC131").Value = 0.8999999999999999

# --- Insert a new fixture before the old row 202, shifting subsequent rows down ---
$ws.Rows.Item(202).Insert()

# --- Row 202: brand-new fixture inserted ---
$ws.Range("A202").Value = 200
$ws.Range("B202").Value = 7623922
$ws.Range("C202").Value = "Costa Rica Primera Division"
$ws.Range("D202").Value = "Costa Rica Primera Division"
$ws.Range("E202").Value = 45353.95833333334
$ws.Range("F202").Value = "Municipal Perez Zeledon"
$ws.Range("G202").Value = "AD Guanacasteca"
$ws.Range("H202").Value = 0
$ws.Range("I202").Value = 0
$ws.Range("J202").Value = "D"
$ws.Range("K202").Value = 2.375
$ws.Range("L202").Value = 3.4
$ws.Range("M202").Value = 2.6
$ws.Range("N202").Value = 2.1
$ws.Range("O202").Value = 3.3
$ws.Range("P202").Value = 3.1
$ws.Range("Q202").Value = -0.25
$ws.Range("R202").Value = 1.875
$ws.Range("S202").Value = 1.925
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 1.825
$ws.Range("V202").Value = 1.975
$ws.Range("W202").Value = -1
$ws.Range("X202").Value = 2.3
$ws.Range("Y202").Value = -1
$ws.Range("Z202").Value = -0.5
$ws.Range("AA202").Value = 0.4625
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = 0.9750000000000001

# --- Row 203 (was old row 202): final score + closing odds now available ---
$ws.Range("A203").Value = 201
$ws.Range("H203").Value = 2
$ws.Range("I203").Value = 2
$ws.Range("J203").Value = "D"
$ws.Range("N203").Value = 1.7
$ws.Range("O203").Value = 3.5
$ws.Range("P203").Value = 4.333
$ws.Range("Q203").Value = -0.75
$ws.Range("R203").Value = 1.95
$ws.Range("S203").Value = 1.85
$ws.Range("W203").Value = -1
$ws.Range("X203").Value = 2.5
$ws.Range("Y203").Value = -1
$ws.Range("Z203").Value = -1
$ws.Range("AA203").Value = 0.8500000000000001
$ws.Range("AB203").Value = 0.8500000000000001
$ws.Range("AC203").Value = -1

# --- Row 204 (was old row 203): final score + closing odds now available ---
$ws.Range("A204").Value = 202
$ws.Range("H204").Value = 2
$ws.Range("I204").Value = 0
$ws.Range("J204").Value = "H"
$ws.Range("N204").Value = 2.55
$ws.Range("P204").Value = 2.5
$ws.Range("R204").Value = 1.925
$ws.Range("S204").Value = 1.875
$ws.Range("W204").Value = 1.55
$ws.Range("X204").Value = -1
$ws.Range("Y204").Value = -1
$ws.Range("Z204").Value = 0.925
$ws.Range("AA204").Value = -1
$ws.Range("AB204").Value = -1
$ws.Range("AC204").Value = 0.8999999999999999

# --- Rows 205-210: brand-new upcoming fixtures appended ---
# Row 205
$ws.Range("A205").Value = 203
$ws.Range("B205").Value = 7623926
$ws.Range("C205").Value = "Costa Rica Primera Division"
$ws.Range("D205").Value = "Costa Rica Primera Division"
$ws.Range("E205").Value = 45359.95833333334
$ws.Range("F205").Value = "Municipal Liberia"
$ws.Range("G205").Value = "Municipal Perez Zeledon"
$ws.Range("K205").Value = 1.666
$ws.Range("L205").Value = 3.8
$ws.Range("M205").Value = 4.5
$ws.Range("N205").Value = 1.666
$ws.Range("O205").Value = 3.8
$ws.Range("P205").Value = 4.75
$ws.Range("Q205").Value = -0.75
$ws.Range("R205").Value = 1.825
$ws.Range("S205").Value = 1.975
$ws.Range("T205").Value = 2.75
$ws.Range("U205").Value = 1.9
$ws.Range("V205").Value = 1.9
$ws.Range("W205").Value = 0
$ws.Range("X205").Value = 0
$ws.Range("Y205").Value = 0
$ws.Range("Z205").Value = 0
$ws.Range("AA205").Value = 0

# Row 206
$ws.Range("A206").Value = 204
$ws.Range("B206").Value = 7920639
$ws.Range("C206").Value = "Costa Rica Primera Division"
$ws.Range("D206").Value = "Costa Rica Primera Division"
$ws.Range("E206").Value = 45360.83333333334
$ws.Range("F206").Value = "Cartagines"
$ws.Range("G206").Value = "Puntarenas"
$ws.Range("K206").Value = 1.666
$ws.Range("L206").Value = 3.75
$ws.Range("M206").Value = 5
$ws.Range("N206").Value = 1.666
$ws.Range("O206").Value = 3.75
$ws.Range("P206").Value = 5
$ws.Range("Q206").Value = -0.75
$ws.Range("R206").Value = 1.85
$ws.Range("S206").Value = 1.95
$ws.Range("T206").Value = 2.75
$ws.Range("U206").Value = 1.95
$ws.Range("V206").Value = 1.85
$ws.Range("W206").Value = 0
$ws.Range("X206").Value = 0
$ws.Range("Y206").Value = 0
$ws.Range("Z206").Value = 0
$ws.Range("AA206").Value = 0

# Row 207
$ws.Range("A207").Value = 205
$ws.Range("B207").Value = 7623923
$ws.Range("C207").Value = "Costa Rica Primera Division"
$ws.Range("D207").Value = "Costa Rica Primera Division"
$ws.Range("E207").Value = 45360.91666666666
$ws.Range("F207").Value = "AD San Carlos"
$ws.Range("G207").Value = "AD Grecia"
$ws.Range("K207").Value = 1.363
$ws.Range("L207").Value = 4.5
$ws.Range("M207").Value = 10
$ws.Range("N207").Value = 1.333
$ws.Range("O207").Value = 4.5
$ws.Range("P207").Value = 10
$ws.Range("Q207").Value = -1.25
$ws.Range("R207").Value = 1.775
$ws.Range("S207").Value = 2.025
$ws.Range("T207").Value = 2.5
$ws.Range("U207").Value = 1.8
$ws.Range("V207").Value = 2
$ws.Range("W207").Value = 0
$ws.Range("X207").Value = 0
$ws.Range("Y207").Value = 0
$ws.Range("Z207").Value = 0
$ws.Range("AA207").Value = 0

# Row 208
$ws.Range("A208").Value = 206
$ws.Range("B208").Value = 7623990
$ws.Range("C208").Value = "Costa Rica Primera Division"
$ws.Range("D208").Value = "Costa Rica Primera Division"
$ws.Range("E208").Value = 45360.95833333334
$ws.Range("F208").Value = "Herediano"
$ws.Range("G208").Value = "Sporting San Jose"
$ws.Range("K208").Value = 1.444
$ws.Range("L208").Value = 4.333
$ws.Range("M208").Value = 7.5
$ws.Range("N208").Value = 1.444
$ws.Range("O208").Value = 4.333
$ws.Range("P208").Value = 7.5
$ws.Range("Q208").Value = -1.25
$ws.Range("R208").Value = 1.975
$ws.Range("S208").Value = 1.825
$ws.Range("T208").Value = 2.5
$ws.Range("U208").Value = 1.825
$ws.Range("V208").Value = 1.975
$ws.Range("W208").Value = 0
$ws.Range("X208").Value = 0
$ws.Range("Y208").Value = 0
$ws.Range("Z208").Value = 0
$ws.Range("AA208").Value = 0

# Row 209
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 7623925
$ws.Range("C209").Value = "Costa Rica Primera Division"
$ws.Range("D209").Value = "Costa Rica Primera Division"
$ws.Range("E209").Value = 45361.75
$ws.Range("F209").Value = "AD Guanacasteca"
$ws.Range("G209").Value = "Alajuelense"
$ws.Range("K209").Value = 5
$ws.Range("L209").Value = 3.5
$ws.Range("M209").Value = 1.727
$ws.Range("N209").Value = 5
$ws.Range("O209").Value = 3.5
$ws.Range("P209").Value = 1.727
$ws.Range("Q209").Value = 0.75
$ws.Range("R209").Value = 1.85
$ws.Range("S209").Value = 1.95
$ws.Range("T209").Value = 2.5
$ws.Range("U209").Value = 2.025
$ws.Range("V209").Value = 1.775
$ws.Range("W209").Value = 0
$ws.Range("X209").Value = 0
$ws.Range("Y209").Value = 0
$ws.Range("Z209").Value = 0
$ws.Range("AA209").Value = 0

# Row 210
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 7916002
$ws.Range("C210").Value = "Costa Rica Primera Division"
$ws.Range("D210").Value = "Costa Rica Primera Division"
$ws.Range("E210").Value = 45361.79166666666
$ws.Range("F210").Value = "Santos de Gupiles"
$ws.Range("G210").Value = "Deportivo Saprissa"
$ws.Range("K210").Value = 6.5
$ws.Range("L210").Value = 4.75
$ws.Range("M210").Value = 1.4
$ws.Range("N210").Value = 6
$ws.Range("O210").Value = 4.333
$ws.Range("P210").Value = 1.45
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = 1.975
$ws.Range("S210").Value = 1.825
$ws.Range("T210").Value = 2.5
$ws.Range("U210").Value = 1.8
$ws.Range("V210").Value = 2
$ws.Range("W210").Value = 0
$ws.Range("X210").Value = 0
$ws.Range("Y210").Value = 0
$ws.Range("Z210").Value = 0
$ws.Range("AA210").Value = 0

